$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.ChartObjects(1).Delete()
$shp = $ws.Shapes.AddChart2(-1, -4169)
$chart = $shp.Chart
$chart.ChartType = -4169
$chart.SetSourceData($ws.Range("B3:C20"), 2)
